# Client: Add MoveState Stamina Control System
# 스테미나가 없으면 달리지 못하게 수정
#
# Tripling the movement "Value" speed multipliers on the Move sheet so that
# running/walking/dashing speeds line up with the new stamina-gated move
# system (Ch_Walk / Ch_Run / Ch_Dash / Ch_HoldWalk / Ch_HoldRun).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Move")

$ws.Range("C3").Value = 0.6   # Ch_Walk:    0.2 -> 0.6
$ws.Range("C4").Value = 1.2   # Ch_Run:     0.4 -> 1.2
$ws.Range("C6").Value = 1.8   # Ch_Dash:    0.6 -> 1.8
$ws.Range("C7").Value = 0.4   # Ch_HoldWalk: 0.17 -> 0.4
$ws.Range("C8").Value = 0.9   # Ch_HoldRun: 0.3 -> 0.9

# The workbook was re-saved with the "Move" tab active/selected.
$ws.Activate()
$ws.Range("D8").Select() | Out-Null
